$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply the "last row" border style (currently on row 71, the final
#        row of the table) to row 46, which will become the new final row
#        once the second worker's rows are removed. Do this before any
#        deletion so the source format on row 71 is still available.
$ws.Range("B71:J71").Copy()
$ws.Range("B46:J46").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Re-sort the first worker's period rows (16-46) into ascending
#        chronological order (currently descending, 2103 down to 1809).
$srcRange = $ws.Range("B16:G46")
$data = $srcRange.Value2

$nRows = $data.GetLength(0)
$nCols = $data.GetLength(1)

$newData = New-Object 'object[,]' $nRows, $nCols
for ($r = 1; $r -le $nRows; $r++) {
    $srcR = $nRows - $r + 1
    for ($c = 1; $c -le $nCols; $c++) {
        $newData[$r - 1, $c - 1] = $data[$srcR, $c]
    }
}

$srcRange.Value2 = $newData

# --- 3. Remove the second worker (LELLYS OROZCO CASTRO) entirely - rows
#        47-71 held that worker's 25 arrears periods.
$ws.Rows("47:71").Delete()

# --- 4. Update the summary figures: total arrears (sum of the remaining
#        worker's periods) and worker count.
$ws.Range("E11").Value = 952053
$ws.Range("C13").Value = 1
